$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the oldest quarter (column D) and shift everything left by one column.
$ws.Range("D1:D59").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)

# The now-empty column M needs the same formatting as column L (which it is a continuation of).
$ws.Range("L1:L59").Copy() | Out-Null
$ws.Range("M1:M59").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

# New column M is the newest quarter: Q4 ending 1401/12, published 1402-02-30.
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"
$ws.Range("M9").Value = "1402-02-30 (2)"

$ws.Range("M12").Value = 5072688
$ws.Range("M13").Value = 0
$ws.Range("M14").Value = 16378616
$ws.Range("M15").Value = 9271779
$ws.Range("M16").Value = 3633218
$ws.Range("M17").Value = 0
$ws.Range("M18").Value = 34356301
$ws.Range("M19").Value = 24896
$ws.Range("M20").Value = 93766
$ws.Range("M21").Value = 0
$ws.Range("M22").Value = 1657549
$ws.Range("M23").Value = 1617
$ws.Range("M24").Value = "-"
$ws.Range("M25").Value = 651
$ws.Range("M26").Value = 1778479
$ws.Range("M27").Value = 36134780
$ws.Range("M29").Value = 3559877
$ws.Range("M30").Value = "-"
$ws.Range("M31").Value = 26602
$ws.Range("M32").Value = 1538046
$ws.Range("M33").Value = 4974162
$ws.Range("M34").Value = 14737294
$ws.Range("M35").Value = 0
$ws.Range("M36").Value = 0
$ws.Range("M37").Value = 24835981
$ws.Range("M38").Value = 0
$ws.Range("M39").Value = "-"
$ws.Range("M40").Value = 0
$ws.Range("M41").Value = 206480
$ws.Range("M42").Value = 206480
$ws.Range("M43").Value = 25042461
$ws.Range("M45").Value = 2250000
$ws.Range("M46").Value = 0
$ws.Range("M47").Value = 0
$ws.Range("M48").Value = -32125
$ws.Range("M49").Value = 86
$ws.Range("M50").Value = 225000
$ws.Range("M51").Value = 0
$ws.Range("M52").Value = "-"
$ws.Range("M53").Value = 0
$ws.Range("M54").Value = "-"
$ws.Range("M55").Value = 0
$ws.Range("M56").Value = 8649358
$ws.Range("M57").Value = 11092319
$ws.Range("M58").Value = 36134780

# Also update the "6th" duplicate-looking publish-date string used for column I (fiscal Q4 1400/12 report).
$ws.Range("I9").Value = "1402-02-30 (8)"

# Column M is now a Q4/year-end column like E and I -- give it the wider 31-width formatting.
# (Excel's ColumnWidth property is offset from the stored OOXML width by ~0.8333 character units.)
$ws.Columns("M").ColumnWidth = 30.166666666666668
